$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "44.374.02"
Set-TextValue $ws.Range("E2") "  +3.56%  "

Set-TextValue $ws.Range("D3") "2.271.40"
Set-TextValue $ws.Range("E3") "  +2.88%  "

Set-TextValue $ws.Range("E4") "  +0.09%  "

Set-TextValue $ws.Range("D5") "322.57"
Set-TextValue $ws.Range("E5") "  +2.39%  "

Set-TextValue $ws.Range("D6") "104.93"
Set-TextValue $ws.Range("E6") "  +6.04%  "

Set-TextValue $ws.Range("D7") "0.590"
Set-TextValue $ws.Range("E7") "  +0.50%  "

Set-TextValue $ws.Range("E8") "  +0.11%  "

Set-TextValue $ws.Range("D9") "0.570"
Set-TextValue $ws.Range("E9") "  +2.20%  "

Set-TextValue $ws.Range("D10") "38.65"
Set-TextValue $ws.Range("E10") "  +5.08%  "

Set-TextValue $ws.Range("D11") "0.0843"
Set-TextValue $ws.Range("E11") "  +2.32%  "

Set-TextValue $ws.Range("D12") "7.89"
Set-TextValue $ws.Range("E12") "  +2.98%  "

Set-TextValue $ws.Range("E13") "  +0.08%  "

Set-TextValue $ws.Range("D14") "0.883"
Set-TextValue $ws.Range("E14") "  +3.22%  "

Set-TextValue $ws.Range("D15") "2.619.94"
Set-TextValue $ws.Range("E15") "  +3.13%  "

Set-TextValue $ws.Range("D16") "14.56"
Set-TextValue $ws.Range("E16") "  +2.61%  "

Set-TextValue $ws.Range("D17") "2.274.02"
Set-TextValue $ws.Range("E17") "  +3.54%  "

Set-TextValue $ws.Range("D18") "44.303.05"
Set-TextValue $ws.Range("E18") "  +3.66%  "

Set-TextValue $ws.Range("D19") "13.85"
Set-TextValue $ws.Range("E19") "  -3.55%  "

Set-TextValue $ws.Range("E20") "  +4.53%  "

Set-TextValue $ws.Range("D21") "6.52"
Set-TextValue $ws.Range("E21") "  +1.74%  "

Set-TextValue $ws.Range("D22") "66.36"
Set-TextValue $ws.Range("E22") "  +1.91%  "

Set-TextValue $ws.Range("E23") "  +2.12%  "

Set-TextValue $ws.Range("D24") "240.28"
Set-TextValue $ws.Range("E24") "  +1.97%  "

Set-TextValue $ws.Range("D25") "2.23"
Set-TextValue $ws.Range("E25") "  +5.04%  "

Set-TextValue $ws.Range("E26") "  +0.20%  "

Set-TextValue $ws.Range("D27") "10.23"
Set-TextValue $ws.Range("E27") "  +2.84%  "

Set-TextValue $ws.Range("D28") "38.40"
Set-TextValue $ws.Range("E28") "  +11.97%  "

Set-TextValue $ws.Range("E29") "  -0.43%  "

Set-TextValue $ws.Range("D30") "6.49"
Set-TextValue $ws.Range("E30") "  +3.37%  "

Set-TextValue $ws.Range("D31") "20.68"
Set-TextValue $ws.Range("E31") "  +1.03%  "

Set-TextValue $ws.Range("D32") "162.52"
Set-TextValue $ws.Range("E32") "  +5.57%  "

Set-TextValue $ws.Range("D33") "0.0882"
Set-TextValue $ws.Range("E33") "  -0.26%  "

Set-TextValue $ws.Range("D34") "2.77"
Set-TextValue $ws.Range("E34") "  -0.21%  "

Set-TextValue $ws.Range("D35") "0.117"
Set-TextValue $ws.Range("E35") "  +9.22%  "

Set-TextValue $ws.Range("E36") "  +5.51%  "

Set-TextValue $ws.Range("D37") "3.13"
Set-TextValue $ws.Range("E37") "  +2.76%  "

Set-TextValue $ws.Range("E38") "  +0.60%  "

Set-TextValue $ws.Range("D39") "3.93"
Set-TextValue $ws.Range("E39") "  +4.55%  "

Set-TextValue $ws.Range("D40") "4.42"
Set-TextValue $ws.Range("E40") "  +0.41%  "

Set-TextValue $ws.Range("D41") "0.0329"
Set-TextValue $ws.Range("E41") "  +1.78%  "

Set-TextValue $ws.Range("D42") "15.53"
Set-TextValue $ws.Range("E42") "  +27.50%  "

Set-TextValue $ws.Range("E43") "  +0.25%  "

Set-TextValue $ws.Range("D44") "1.788.10"
Set-TextValue $ws.Range("E44") "  -2.35%  "

Set-TextValue $ws.Range("D45") "0.209"
Set-TextValue $ws.Range("E45") "  +1.11%  "

Set-TextValue $ws.Range("D46") "86.11"
Set-TextValue $ws.Range("E46") "  -2.85%  "

Set-TextValue $ws.Range("D47") "5.44"
Set-TextValue $ws.Range("E47") "  +1.62%  "

Set-TextValue $ws.Range("D48") "60.63"
Set-TextValue $ws.Range("E48") "  +0.01%  "

Set-TextValue $ws.Range("B49") "ordi"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D49") "75.56"
Set-TextValue $ws.Range("E49") "  +0.55%  "

Set-TextValue $ws.Range("B50") "Stacks"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "1.73"
Set-TextValue $ws.Range("E50") "  +9.70%  "

Set-TextValue $ws.Range("D51") "104.12"
Set-TextValue $ws.Range("E51") "  +1.43%  "
